{"js": "// Translate the Heart Failure risk-model table (caption, headers, footnote)\n// from English to Portuguese.\nconst replacements = [\n  [\"Risk model for Heart Failure\", \"Modelo de Risco para Insufici\u00eancia Card\u00edaca\"],\n  [\"By Age and Sex\", \"Por Idade e Sexo\"],\n  [\"Predictor\", \"Preditor\"],\n  [\"Odds Ratio\", \"Raz\u00e3o de Probabilidades\"],\n  [\"p Value\", \"Valor p\"],\n  [\"Low CI\", \"IC Inferior\"],\n  [\"High CI\", \"IC Superior\"],\n  [\"Source: MlR dataset\", \"Fonte: Conjunto de dados MlR\"],\n];\n\nconst body = context.document.body;\n\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Translate the Heart Failure risk-model table (caption, headers, footnote)\n# from English to Portuguese.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Risk model for Heart Failure\"; Replace = \"Modelo de Risco para Insufici\u00eancia Card\u00edaca\" },\n    @{ Find = \"By Age and Sex\";                Replace = \"Por Idade e Sexo\" },\n    @{ Find = \"Predictor\";                     Replace = \"Preditor\" },\n    @{ Find = \"Odds Ratio\";                    Replace = \"Raz\u00e3o de Probabilidades\" },\n    @{ Find = \"p Value\";                       Replace = \"Valor p\" },\n    @{ Find = \"Low CI\";                        Replace = \"IC Inferior\" },\n    @{ Find = \"High CI\";                       Replace = \"IC Superior\" },\n    @{ Find = \"Source: MlR dataset\";           Replace = \"Fonte: Conjunto de dados MlR\" }\n)\n\nforeach ($item in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $item.Find\n    $find.Replacement.Text = $item.Replace\n    $find.Execute($item.Find, $true, $true, $false, $false, $false, $true, 1, $false, $item.Replace, 2)\n}\n"}
